$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Data change: InstructorPaymentFree (Sheet2 row 19) switches from a
# checkbox/bit field ("bit(1)" / "b'0'," / "x") to an integer/instructor
# payment field ("int(11)" / "'0'," / "i"). Sheet1's VLOOKUP-driven columns
# (E26, and the running-concatenation F26:F33) recalc automatically.
$ws2.Range("C19").Value = "int(11)"

# The target value starts with a literal apostrophe ('0',) which Excel's
# COM Value setter would otherwise treat as a text-prefix quote character
# and strip. Route it through a formula + paste-values round trip so the
# literal apostrophe is preserved as a normal shared-string value.
$ws2.Range("K1").Formula = "=""'0',"""
$ws2.Range("K1").Copy()
$ws2.Range("E19").PasteSpecial(-4163)  # xlPasteValues
$ws2.Range("K1").ClearContents()
$excel.CutCopyMode = $false

$ws2.Range("I19").Value = "i"

# --- View change: Sheet2's window is scrolled so row 18 is the top-left
# visible row, with I20 as the active selected cell.
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("I20").Select()

# Restore Sheet1 as the active/displayed sheet with its original selection
# so its sheetView stays exactly as it was before this edit.
$ws1.Activate()
$ws1.Range("F26").Select()
